# Apply the "gen" sheet renewable-share columns + restore the user's
# navigation/selection state (sheet tab, zoom, active cell) on each
# worksheet, matching the authored commit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "gen" sheet: add Gas/Wind/PV share columns (K, L, M)
# ---------------------------------------------------------------------
$wsGen = $wb.Worksheets.Item("gen")

# Add the new shared strings in the same order they appear in the
# rebuilt sharedStrings table: Wind Share, PV Share, Gas Share.
$wsGen.Range("L1").Value = "Wind Share"
$wsGen.Range("M1").Value = "PV Share"
$wsGen.Range("K1").Value = "Gas Share"

# All 27 DG rows (2-28) get Gas Share = 0, Wind Share = 0, PV Share = 1
for ($r = 2; $r -le 28; $r++) {
    $wsGen.Cells.Item($r, 11).Value = 0
    $wsGen.Cells.Item($r, 12).Value = 0
    $wsGen.Cells.Item($r, 13).Value = 1
}

# Approximate the column widths Excel auto-sized for the new columns.
$wsGen.Range("K1").ColumnWidth = 9.3
$wsGen.Range("L1").ColumnWidth = 9.3

# ---------------------------------------------------------------------
# 2) Restore each sheet's view/selection state
# ---------------------------------------------------------------------

# "bus": scrolled back to top, selection moved to G5
$wsBus = $wb.Worksheets.Item("bus")
$wsBus.Activate() | Out-Null
$wsBus.Range("G5").Select() | Out-Null

# "line&trafo": becomes the active tab, selection moved to K30
$wsLine = $wb.Worksheets.Item("line&trafo")
$wsLine.Activate() | Out-Null
$wsLine.Range("K30").Select() | Out-Null

# "load": selection stays at D32 (just loses the active-tab flag)
$wsLoad = $wb.Worksheets.Item("load")
$wsLoad.Activate() | Out-Null
$wsLoad.Range("D32").Select() | Out-Null

# "gen": zoom changed to 115%, selection moved to E21
$wsGen.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 115
$wsGen.Range("E21").Select() | Out-Null

# Leave "line&trafo" as the final active sheet/tab.
$wsLine.Activate() | Out-Null
